$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (rich-text shared strings) ---
# A8: "Volume 30   Number  2" -> "Volume 30   Number  3"
$volCell = $ws.Range("A8")
$volCell.Characters(21, 1).Text = "3"

# C9: "Report Covering the Week  1/9/2023  Through  1/15/2023"
#  -> "Report Covering the Week  1/16/2023  Through  1/22/2023"
# Replace the later date first so the earlier replacement's position is unaffected.
$weekCell = $ws.Range("C9")
$weekCell.Characters(46, 9).Text = "1/22/2023"
$weekCell.Characters(27, 8).Text = "1/16/2023"

# --- Style donor formats for cells that flip from text placeholder to numeric ---
$numFmt16 = $ws.Range("F22").NumberFormat
$align16 = $ws.Range("F22").HorizontalAlignment
$numFmt15 = $ws.Range("H22").NumberFormat
$align15 = $ws.Range("H22").HorizontalAlignment

# --- Crime Complaints table value updates (rows 15-29) ---
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -50
$ws.Range("N15").Value = 200

$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 73
$ws.Range("H16").Value = 58.695652173913
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 38
$ws.Range("K16").Value = 34.210526315789
$ws.Range("L16").Value = 218.75
$ws.Range("M16").Value = 155
$ws.Range("N16").Value = -63.309352517985

$ws.Range("C17").Value = 20
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 77
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = 60.416666666666
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 91.176470588235
$ws.Range("L17").Value = 71.052631578947
$ws.Range("M17").Value = 209.52380952381
$ws.Range("N17").Value = 20.37037037037

$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 80
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = 83.333333333333
$ws.Range("L18").Value = 57.142857142857
$ws.Range("M18").Value = -29.032258064516
$ws.Range("N18").Value = -76.344086021505

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -6.756756756756
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 59
$ws.Range("K19").Value = -10.169491525423
$ws.Range("L19").Value = 96.296296296296
$ws.Range("M19").Value = 76.666666666666
$ws.Range("N19").Value = 70.967741935483

$ws.Range("C20").Value = 12
$ws.Range("E20").Value = -29.411764705882
$ws.Range("F20").Value = 73
$ws.Range("G20").Value = 56
$ws.Range("H20").Value = 30.357142857142
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = 1.88679245283
$ws.Range("L20").Value = 390.909090909091
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = -62.758620689655

$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = 6.55737704918
$ws.Range("F21").Value = 319
$ws.Range("G21").Value = 251
$ws.Range("H21").Value = 27.091633466135
$ws.Range("I21").Value = 248
$ws.Range("J21").Value = 202
$ws.Range("K21").Value = 22.772277227722
$ws.Range("L21").Value = 131.775700934579
$ws.Range("M21").Value = 111.965811965812
$ws.Range("N21").Value = -46.781115879828

$ws.Range("C22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("M22").Value = 0
$ws.Range("C22").NumberFormat = $numFmt16
$ws.Range("C22").HorizontalAlignment = $align16
$ws.Range("I22").NumberFormat = $numFmt16
$ws.Range("I22").HorizontalAlignment = $align16
$ws.Range("M22").NumberFormat = $numFmt15
$ws.Range("M22").HorizontalAlignment = $align15

$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 28
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 40
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 84.615384615384
$ws.Range("M23").Value = 118.181818181818

$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 57.692307692307
$ws.Range("F24").Value = 149
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = 44.660194174757
$ws.Range("I24").Value = 111
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 38.75
$ws.Range("L24").Value = 56.338028169014
$ws.Range("M24").Value = 26.136363636363

$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 79
$ws.Range("G25").Value = 67
$ws.Range("H25").Value = 17.910447761194
$ws.Range("I25").Value = 69
$ws.Range("J25").Value = 55
$ws.Range("K25").Value = 25.454545454545
$ws.Range("L25").Value = 60.465116279069
$ws.Range("M25").Value = -10.38961038961

$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -66.666666666666
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 5
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -28.571428571428

$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = -40
$ws.Range("C27").NumberFormat = $numFmt16
$ws.Range("C27").HorizontalAlignment = $align16
$ws.Range("D27").NumberFormat = $numFmt16
$ws.Range("D27").HorizontalAlignment = $align16
$ws.Range("E27").NumberFormat = $numFmt15
$ws.Range("E27").HorizontalAlignment = $align15

$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 1
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -66.666666666666
$ws.Range("N28").Value = -90.90909090909
$ws.Range("C28").NumberFormat = $numFmt16
$ws.Range("C28").HorizontalAlignment = $align16
$ws.Range("I28").NumberFormat = $numFmt16
$ws.Range("I28").HorizontalAlignment = $align16

$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 1
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -66.666666666666
$ws.Range("N29").Value = -90
$ws.Range("C29").NumberFormat = $numFmt16
$ws.Range("C29").HorizontalAlignment = $align16
$ws.Range("I29").NumberFormat = $numFmt16
$ws.Range("I29").HorizontalAlignment = $align16

